$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Jalon 2")

# Update durée idéale (column D) values for rows 3-5
$ws.Range("D3").Value = 6
$ws.Range("D4").Value = 20
$ws.Range("D5").Value = 8

# Add new activity rows (13-16 in the "Numéro" column, data rows 14-17)
$ws.Range("A14").Value = 13
$ws.Range("B14").Value = "menu responsive"

$ws.Range("A15").Value = 14
$ws.Range("B15").Value = "menu responsive"

$ws.Range("A16").Value = 15
$ws.Range("B16").Value = "menu responsive"

$ws.Range("A17").Value = 16
$ws.Range("B17").Value = "Galerie responsive sur page d'accueil"

# Move the active selection to D6, matching the saved workbook state
$ws.Range("D6").Select() | Out-Null
